$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 2.37

# Row 4 updates
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 4.2
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("S4").Value = 2.15
$ws.Range("T4").Value = 1.67
$ws.Range("Y4").Value = 1.44
$ws.Range("Z4").Value = 2.63
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 7.5
$ws.Range("AF4").Value = 13
$ws.Range("AK4").Value = 19
$ws.Range("AN4").Value = 10
$ws.Range("AO4").Value = 21
$ws.Range("AP4").Value = 15
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 41

# Row 5 updates
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
